# Pedalboard HW BoM refresh — reference designator renumbering.
# Applies the "References" column (column D) updates on the BoM and DNF
# sheets to match the regenerated KiCad BoM (components renumbered).

$wb  = $excel.ActiveWorkbook
$bom = $wb.Worksheets.Item("BoM")
$dnf = $wb.Worksheets.Item("DNF")

# --- BoM sheet -------------------------------------------------------
$bom.Range("D9").Value  = "C22 C26"
$bom.Range("D10").Value = "C29 C34"
$bom.Range("D12").Value = "C1 C11 C13 C15 C16 C17 C18 C19 C20 C21 C23 C25 C30 C35"
$bom.Range("D13").Value = "C12 C14"
$bom.Range("D14").Value = "C2 C3 C4 C5 C6 C7 C8 C9 C10 C24 C27 C28 C32 C36 C37"
$bom.Range("D15").Value = "C33"

$bom.Range("D17").Value = "D7"
$bom.Range("D18").Value = "D6"
$bom.Range("D19").Value = "D2 D3"
$bom.Range("D21").Value = "D4 D5"

$bom.Range("D24").Value = "J27"
$bom.Range("D26").Value = "J5 J8 J18 J19 J20 J22"
$bom.Range("D27").Value = "J9 J10 J13 J15 J17 J21 J23 J24"
$bom.Range("D28").Value = "J14"
$bom.Range("D29").Value = "J28"
$bom.Range("D30").Value = "J11"

$bom.Range("D34").Value = "R5 R13"
$bom.Range("D35").Value = "R6 R10"
$bom.Range("D36").Value = "R4"
$bom.Range("D37").Value = "R3 R7"
$bom.Range("D38").Value = "R1 R2 R8 R11 R16"
$bom.Range("D39").Value = "R15 R17"
$bom.Range("D40").Value = "R12"
$bom.Range("D41").Value = "R14"

$bom.Range("D43").Value = "SW5 SW6 SW7 SW8 SW9 SW10"

$bom.Range("D44").Value = "U1 U7"
$bom.Range("D47").Value = "U3"
$bom.Range("D48").Value = "U2"
$bom.Range("D49").Value = "U4"

# --- DNF sheet --------------------------------------------------------
$dnf.Range("D10").Value = "J12 J16"
$dnf.Range("D11").Value = "J25 J26"
$dnf.Range("D13").Value = "R9"
$dnf.Range("D14").Value = "SW3 SW4"
